$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 143, shifting existing row 143 (and below) down to 144.
$ws.Rows.Item(143).Insert()

# Populate the newly inserted row 143 with the new data record.
$ws.Range("A143").Value = 4
$ws.Range("B143").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C143").Value = "Los Lagos"
$ws.Range("D143").Value = 44606
$ws.Range("E143").Value = 10
$ws.Range("F143").Value = 100112043
$ws.Range("G143").Value = "Pepino ensalada"
$ws.Range("H143").Value = "Sin especificar"
$ws.Range("I143").Value = "Primera"
$ws.Range("J143").Value = 200
$ws.Range("K143").Value = 17000
$ws.Range("L143").Value = 17000
$ws.Range("M143").Value = 17000
$ws.Range("N143").Value = "$/caja 70 unidades"
$ws.Range("O143").Value = "Región del Maule"
$ws.Range("P143").Value = 243
$ws.Range("Q143").Value = 70
$ws.Range("R143").Value = "Hortaliza"
